$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$strMap = @{
    20 = "ECs"
    21 = "FAPs"
    22 = "MuSCs"
    23 = "Col4a1"
    24 = "Itgb8"
}

$rows = @(
    @{ A=$strMap[20]; B=$strMap[23]; C=$strMap[24]; D=$strMap[20]; E=3; F=1; G=93.39526366666666; H=280.185791; I=0.2167755775732346; J=0.2167755775732346; K=1; L=0.3333333333333333; M=0.027767; N=0.083301; O=0.002923627791763407; P=0.002923627791763407; Q=2.593306286232334; R=23.339756576091; S=0.000633771103168673; T=0.0006337711031686729 },
    @{ A=$strMap[20]; B=$strMap[23]; C=$strMap[24]; D=$strMap[21]; E=3; F=1; G=93.39526366666666; H=280.185791; I=0.2167755775732346; J=0.2167755775732346; K=3; L=1; M=4.237840333333334; N=12.713521; O=0.4462083687682946; P=0.4462083687682946; Q=395.7942153089012; R=3562.147937780111; S=0.09672707685775792; T=0.09672707685775789 },
    @{ A=$strMap[20]; B=$strMap[23]; C=$strMap[24]; D=$strMap[22]; E=3; F=1; G=93.39526366666666; H=280.185791; I=0.2167755775732346; J=0.2167755775732346; K=3; L=1; M=5.231839666666667; N=15.695519; O=0.5508680034399419; P=0.5508680034399419; Q=488.6290451300588; R=4397.66140617053; S=0.119414729612308; T=0.119414729612308 },
    @{ A=$strMap[21]; B=$strMap[23]; C=$strMap[24]; D=$strMap[20]; E=3; F=1; G=239.807332; H=719.421996; I=0.5566060939249745; J=0.5566060939249745; K=1; L=0.3333333333333333; M=0.027767; N=0.083301; O=0.002923627791763407; P=0.002923627791763407; Q=6.658730187644; R=59.928571688796; S=0.001627309045263929; T=0.001627309045263929 },
    @{ A=$strMap[21]; B=$strMap[23]; C=$strMap[24]; D=$strMap[21]; E=3; F=1; G=239.807332; H=719.421996; I=0.5566060939249745; J=0.5566060939249745; K=3; L=1; M=4.237840333333334; N=12.713521; O=0.4462083687682946; P=0.4462083687682946; Q=1016.265183778657; R=9146.386654007916; S=0.2483622972167551; T=0.248362297216755 },
    @{ A=$strMap[21]; B=$strMap[23]; C=$strMap[24]; D=$strMap[22]; E=3; F=1; G=239.807332; H=719.421996; I=0.5566060939249745; J=0.5566060939249745; K=3; L=1; M=5.231839666666667; N=15.695519; O=0.5508680034399419; P=0.5508680034399419; Q=1254.633511915103; R=11291.70160723593; S=0.3066164876629555; T=0.3066164876629555 },
    @{ A=$strMap[22]; B=$strMap[23]; C=$strMap[24]; D=$strMap[20]; E=3; F=1; G=97.63589966666666; H=292.907699; I=0.226618328501791; J=0.2266183285017909; K=1; L=0.3333333333333333; M=0.027767; N=0.083301; O=0.002923627791763407; P=0.002923627791763407; Q=2.711056026044333; R=24.399504234399; S=0.0006625476433308056; T=0.0006625476433308055 },
    @{ A=$strMap[22]; B=$strMap[23]; C=$strMap[24]; D=$strMap[21]; E=3; F=1; G=97.63589966666666; H=292.907699; I=0.226618328501791; J=0.2266183285017909; K=3; L=1; M=4.237840333333334; N=12.713521; O=0.4462083687682946; P=0.4462083687682946; Q=413.7653535886865; R=3723.888182298179; S=0.1011189946937817; T=0.1011189946937817 },
    @{ A=$strMap[22]; B=$strMap[23]; C=$strMap[24]; D=$strMap[22]; E=3; F=1; G=97.63589966666666; H=292.907699; I=0.226618328501791; J=0.2266183285017909; K=3; L=1; M=5.231839666666667; N=15.695519; O=0.5508680034399419; P=0.5508680034399419; Q=510.8153727667534; R=4597.338354900781; S=0.1248367861646785; T=0.1248367861646785 },
)

$rowIndex = 2
foreach ($rowData in $rows) {
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowIndex").Value = $rowData[$col]
    }
    $rowIndex++
}

Write-Output "done"